# Apply the "Updated cryptos list" price-refresh edit.
# A new coin (WrappedliquidstakedEther2.0) was inserted at row 21, shifting
# rows 21-50 down by one (NEARProtocol drops off the bottom of the list),
# and every Price / Volume(1h) figure is refreshed to the latest snapshot.
# Price figures are stored as literal text (e.g. "4.100", "1.847.97") in the
# source data, so force Text number-format before writing them to stop Excel
# from "helpfully" re-parsing them as numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.177.39"
$ws.Cells.Item(2, 5).Value = "  -1.71%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.847.97"
$ws.Cells.Item(3, 5).Value = "  -2.62%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  +0.14%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "232.79"

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.13%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4684"
$ws.Cells.Item(7, 5).Value = "  -2.65%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2708"
$ws.Cells.Item(8, 5).Value = "  -4.66%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06365"
$ws.Cells.Item(9, 5).Value = "  -2.84%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.819.86"
$ws.Cells.Item(10, 5).Value = "  -5.49%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07424"
$ws.Cells.Item(11, 5).Value = "  -0.81%  "

$ws.Cells.Item(12, 5).Value = "  -2.70%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.944"
$ws.Cells.Item(13, 5).Value = "  -3.29%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "84.93"
$ws.Cells.Item(14, 5).Value = "  -3.68%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.6258"
$ws.Cells.Item(15, 5).Value = "  -6.44%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "30.120.80"
$ws.Cells.Item(16, 5).Value = "  -1.83%  "

$ws.Cells.Item(17, 5).Value = "  +0.13%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "229.08"
$ws.Cells.Item(18, 5).Value = "  -0.07%  "

$ws.Cells.Item(19, 5).Value = "  -5.33%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.000007298"
$ws.Cells.Item(20, 5).Value = "  -4.27%  "

$ws.Cells.Item(21, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "2.101.12"
$ws.Cells.Item(21, 5).Value = "  -3.84%  "

$ws.Cells.Item(22, 2).Value = "BinanceUSD"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.002"
$ws.Cells.Item(22, 5).Value = "  +0.10%  "

$ws.Cells.Item(23, 2).Value = "Uniswap"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.938"
$ws.Cells.Item(23, 5).Value = "  -7.08%  "

$ws.Cells.Item(24, 2).Value = "Chainlink"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.928"
$ws.Cells.Item(24, 5).Value = "  -4.74%  "

$ws.Cells.Item(25, 2).Value = "Cosmos"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.226"
$ws.Cells.Item(25, 5).Value = "  -0.54%  "

$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "166.06"
$ws.Cells.Item(26, 5).Value = "  -2.38%  "

$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "17.75"
$ws.Cells.Item(27, 5).Value = "  -4.98%  "

$ws.Cells.Item(28, 2).Value = "LidoDAOToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.863"
$ws.Cells.Item(28, 5).Value = "  -4.72%  "

$ws.Cells.Item(29, 2).Value = "Stellar"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.1036"
$ws.Cells.Item(29, 5).Value = "  +2.62%  "

$ws.Cells.Item(30, 2).Value = "Toncoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.389"
$ws.Cells.Item(30, 5).Value = "  -1.00%  "

$ws.Cells.Item(31, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "4.100"
$ws.Cells.Item(31, 5).Value = "  -6.13%  "

$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.865"
$ws.Cells.Item(32, 5).Value = "  -4.17%  "

$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.04875"
$ws.Cells.Item(33, 5).Value = "  -3.91%  "

$ws.Cells.Item(34, 2).Value = "ARBITRUM"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.157"
$ws.Cells.Item(34, 5).Value = "  -5.55%  "

$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.7116"
$ws.Cells.Item(35, 5).Value = "  -5.57%  "

$ws.Cells.Item(36, 2).Value = "Frax"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.000"
$ws.Cells.Item(36, 5).Value = "  -0.74%  "

$ws.Cells.Item(37, 2).Value = "HuobiToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.698"
$ws.Cells.Item(37, 5).Value = "  -0.52%  "

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.01852"
$ws.Cells.Item(38, 5).Value = "  -1.85%  "

$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.637"
$ws.Cells.Item(39, 5).Value = "  -0.77%  "

$ws.Cells.Item(40, 2).Value = "TrustWalletToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.9081"
$ws.Cells.Item(40, 5).Value = "  -1.18%  "

$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.936"
$ws.Cells.Item(41, 5).Value = "  -7.09%  "

$ws.Cells.Item(42, 2).Value = "Quant"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "105.21"
$ws.Cells.Item(42, 5).Value = "  -1.78%  "

$ws.Cells.Item(43, 2).Value = "PaxDollar"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.9983"
$ws.Cells.Item(43, 5).Value = "  -0.58%  "

$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.547"
$ws.Cells.Item(44, 5).Value = "  -5.29%  "

$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.4061"
$ws.Cells.Item(45, 5).Value = "  -5.53%  "

$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "7.017"
$ws.Cells.Item(46, 5).Value = "  -5.51%  "

$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "60.14"
$ws.Cells.Item(47, 5).Value = "  -7.02%  "

$ws.Cells.Item(48, 2).Value = "Algorand"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.1184"
$ws.Cells.Item(48, 5).Value = "  -7.07%  "

$ws.Cells.Item(49, 2).Value = "EnergySwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "8.635"
$ws.Cells.Item(49, 5).Value = "  -4.59%  "

$ws.Cells.Item(50, 2).Value = "Elrond"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "32.99"
$ws.Cells.Item(50, 5).Value = "  -2.96%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.05572"
$ws.Cells.Item(51, 5).Value = "  -1.63%  "

Write-Output "Updated crypto price rows 2-51"
